# The workbook has two worksheets:
#   1) "Direction 0 STOPS" (currently the selected/active tab)
#   2) "Direction 1 STOPS"
#
# This edit:
#   - clears the (erroneous) Stop ID values in B2:B3 on "Direction 1 STOPS"
#     (keeping their existing number formatting/style intact), so the
#     validator can surface them as missing data instead of stopping after
#     the first discovered row,
#   - makes "Direction 1 STOPS" the active sheet/tab, with B3 selected.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item(2)

# Clear out the B2 and B3 values on "Direction 1 STOPS", leaving formatting in place.
$ws2.Range("B2:B3").ClearContents()

# Switch the active sheet to "Direction 1 STOPS" and select cell B3 there.
$ws2.Activate()
$ws2.Range("B3").Select()
